$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.803.34'
$ws.Range("E2").Value = '  +1.50%  '
$ws.Range("D3").Value = '2.623.28'
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.27'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.88'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.90%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.70%  '
$ws.Range("D9").Value = '2.622.86'
$ws.Range("E9").Value = '  +1.06%  '
$ws.Range("E10").Value = '  +10.18%  '
$ws.Range("E11").Value = '  +0.86%  '
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.16'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.39%  '
$ws.Range("E15").Value = '  +3.75%  '
$ws.Range("D16").Value = '3.102.06'
$ws.Range("E16").Value = '  +1.64%  '
$ws.Range("D17").Value = '67.765.20'
$ws.Range("E17").Value = '  +1.92%  '
$ws.Range("D18").Value = '2.631.00'
$ws.Range("E18").Value = '  +1.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.35'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '365.18'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.50%  '
$ws.Range("E21").Value = '  -2.52%  '
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("E23").Value = '  +4.70%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.25'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.15'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.76%  '
$ws.Range("E27").Value = '  +3.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '587.16'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.57%  '
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.87'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.132'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.48'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '155.65'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.07%  '
$ws.Range("E40").Value = '  +1.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.44'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  +3.31%  '
$ws.Range("E43").Value = '  +2.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.15'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.44'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '157.48'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.69%  '
$ws.Range("D48").Value = '0.0₆0289'
$ws.Range("E48").Value = '  -7.34%  '
$ws.Range("E49").Value = '  +0.35%  '
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.626'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.26%  '
